$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("I2").Value = 9.4
$ws.Range("J2").Value = 3.8
$ws.Range("P2").Value = 1.8
$ws.Range("U2").Value = 1.73
$ws.Range("V2").Value = 1.12

# Row 4
$ws.Range("F4").Value = 1.77
$ws.Range("G4").Value = 1.89
$ws.Range("H4").Value = 5.4
$ws.Range("I4").Value = 6.6
$ws.Range("K4").Value = 3.85
$ws.Range("N4").Value = 2.82
$ws.Range("P4").Value = 1.62
$ws.Range("Q4").Value = 2.26
$ws.Range("R4").Value = 1.23
$ws.Range("T4").Value = 2.1
$ws.Range("W4").Value = 2.12
$ws.Range("AA4").Value = 1000
$ws.Range("AF4").Value = 12

# Row 5
$ws.Range("G5").Value = 8
$ws.Range("I5").Value = 1.53
$ws.Range("K5").Value = 4.7
$ws.Range("L5").Value = 1.37
$ws.Range("R5").Value = 1.46
$ws.Range("S5").Value = 3.05
$ws.Range("U5").Value = 1.97
$ws.Range("V5").Value = 2.92
$ws.Range("AB5").Value = 26

# Row 6
$ws.Range("H6").Value = 13
$ws.Range("I6").Value = 14
$ws.Range("J6").Value = 7.4
$ws.Range("K6").Value = 7.8
$ws.Range("V6").Value = 1.07
$ws.Range("Z6").Value = 150
$ws.Range("AF6").Value = 10
$ws.Range("AI6").Value = 130
$ws.Range("AK6").Value = 12
$ws.Range("AL6").Value = 27

# Row 7
$ws.Range("F7").Value = 3.4
$ws.Range("G7").Value = 3.45
$ws.Range("H7").Value = 2.26
$ws.Range("N7").Value = 4.8
$ws.Range("O7").Value = 1.24
$ws.Range("P7").Value = 2.26
$ws.Range("U7").Value = 2.48
$ws.Range("V7").Value = 1.78

# Row 8
$ws.Range("I8").Value = 8.4
$ws.Range("R8").Value = 1.44
$ws.Range("AE8").Value = 130

# Row 9
$ws.Range("F9").Value = 3.1
$ws.Range("G9").Value = 3.15
$ws.Range("H9").Value = 2.32
$ws.Range("I9").Value = 2.34
$ws.Range("N9").Value = 6
$ws.Range("O9").Value = 1.19
$ws.Range("P9").Value = 2.64
$ws.Range("Q9").Value = 1.57
$ws.Range("R9").Value = 1.67
$ws.Range("S9").Value = 2.42
$ws.Range("T9").Value = 1.52
$ws.Range("U9").Value = 2.76
$ws.Range("V9").Value = 1.75
$ws.Range("W9").Value = 1.46
$ws.Range("Y9").Value = 16
$ws.Range("Z9").Value = 18
$ws.Range("AA9").Value = 30
$ws.Range("AE9").Value = 21
$ws.Range("AF9").Value = 25
$ws.Range("AH9").Value = 14
$ws.Range("AJ9").Value = 50
$ws.Range("AK9").Value = 29
$ws.Range("AM9").Value = 55
$ws.Range("AN9").Value = 19
$ws.Range("AO9").Value = 11.5

# Row 10
$ws.Range("F10").Value = 2.28
$ws.Range("G10").Value = 2.32
$ws.Range("N10").Value = 5.9
$ws.Range("P10").Value = 2.66
$ws.Range("Q10").Value = 1.57
$ws.Range("R10").Value = 1.68
$ws.Range("V10").Value = 1.45
$ws.Range("X10").Value = 24
$ws.Range("Y10").Value = 19.5
$ws.Range("AO10").Value = 19

# Row 11
$ws.Range("F11").Value = 2.16
$ws.Range("G11").Value = 2.18
$ws.Range("H11").Value = 3.6
$ws.Range("I11").Value = 3.65
$ws.Range("Q11").Value = 1.72
$ws.Range("S11").Value = 2.78
$ws.Range("T11").Value = 1.63
$ws.Range("V11").Value = 1.37
$ws.Range("W11").Value = 1.84
$ws.Range("AA11").Value = 65
$ws.Range("AB11").Value = 12.5
$ws.Range("AC11").Value = 8.6
$ws.Range("AD11").Value = 14.5
$ws.Range("AK11").Value = 20
$ws.Range("AO11").Value = 28

# Row 12
$ws.Range("I12").Value = 19.5
$ws.Range("N12").Value = 10.5
$ws.Range("Q12").Value = 1.26
$ws.Range("R12").Value = 2.42
$ws.Range("S12").Value = 1.65
$ws.Range("T12").Value = 1.81
$ws.Range("U12").Value = 2.16
$ws.Range("X12").Value = 1000
$ws.Range("Z12").Value = 1000
$ws.Range("AB12").Value = 20
$ws.Range("AC12").Value = 29
$ws.Range("AD12").Value = 1000
$ws.Range("AE12").Value = 1000
$ws.Range("AF12").Value = 12
$ws.Range("AG12").Value = 14.5
$ws.Range("AH12").Value = 42
$ws.Range("AJ12").Value = 12.5
$ws.Range("AN12").Value = 2.46

# Row 13
$ws.Range("U13").Value = 2.54
